$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Old order: H=Group, I=Block, J=Room
# New order: H=Room, I=Block, J=Group, K=Level, L=Course
$ws.Range("H1").Value = "Room"
$ws.Range("I1").Value = "Block"
$ws.Range("J1").Value = "Group"
$ws.Range("K1").Value = "Level"
$ws.Range("L1").Value = "Course"

# --- Data rows ---
# Columns: A Day, B Time, C Module Code, D Module Title, E Hours, F Class Type,
#          G Lecturer, H Room, I Block, J Group, K Level, L Course

$data = @(
    @{Row=2;  Day="SUN"; Time="12:30-15:30"; Code="5CS024"; Title="Collaborative Development"; Hours=2.5; Type="Workshop"; Lecturer="Mr. Udaya Kandel";  Room="SR-02 Bilston";  Block="WLV"; Group="L5CG8";          Level=5; Course="BCS"}
    @{Row=3;  Day="TUE"; Time="7:00-9:00";   Code="5CS022"; Title="Human Computer Interaction"; Hours=2;   Type="Lecture";  Lecturer="Mr. Apurba Neupane"; Room="LT-02 Telford";  Block="WLV"; Group="L5CG(5+6+7+8)"; Level=5; Course="BCS"}
    @{Row=4;  Day="TUE"; Time="9:30-11:30";  Code="5CS020"; Title="Distributed and Cloud Systems Programming"; Hours=2; Type="Lecture";  Lecturer="Mr. Sumanta Silwal"; Room="LT-01 Wulfruna"; Block="WLV"; Group="L5CG(5+6+7+8)"; Level=5; Course="BCS"}
    @{Row=5;  Day="WED"; Time="7:00-9:00";   Code="5CS024"; Title="Collaborative Development"; Hours=2;   Type="Lecture";  Lecturer="Mr. Raj Shrestha";   Room="LT-02 Telford";  Block="WLV"; Group="L5CG(5+6+7+8)"; Level=5; Course="BCS"}
    @{Row=6;  Day="WED"; Time="9:30-11:30";  Code="5CS020"; Title="Distributed and Cloud Systems Programming"; Hours=2; Type="Tutorial"; Lecturer="Mr. Prabin Sapkota"; Room="TR-02 Stafford"; Block="WLV"; Group="L5CG8"; Level=5; Course="BCS"}
    @{Row=7;  Day="THU"; Time="7:00-9:00";   Code="5CS022"; Title="Human Computer Interaction"; Hours=2;   Type="Tutorial"; Lecturer="Mr. Apurba Neupane"; Room="TR-02 Stafford"; Block="WLV"; Group="L5CG8"; Level=5; Course="BCS"}
    @{Row=8;  Day="THU"; Time="9:30-12:00";  Code="5CS020"; Title="Distributed and Cloud Systems Programming"; Hours=2.5; Type="Workshop"; Lecturer="Mr. Prabin Sapkota"; Room="Lab-01 Mander"; Block="WLV"; Group="L5CG8"; Level=5; Course="BCS"}
    @{Row=9;  Day="FRI"; Time="7:00-9:30";   Code="5CS022"; Title="Human Computer Interaction"; Hours=2.5; Type="Workshop"; Lecturer="Mr. Apurba Neupane"; Room="SR-04 Crompton"; Block="WLV"; Group="L5CG8"; Level=5; Course="BCS"}
    @{Row=10; Day="FRI"; Time="10:00-12:00"; Code="5CS024"; Title="Collaborative Development"; Hours=2;   Type="Tutorial"; Lecturer="Mr. Udaya Kandel";   Room="SR-03 Wolves";  Block="WLV"; Group="L5CG8"; Level=5; Course="BCS"}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Range("A$r").Value = $rowData.Day
    $ws.Range("B$r").Value = $rowData.Time
    $ws.Range("C$r").Value = $rowData.Code
    $ws.Range("D$r").Value = $rowData.Title
    $ws.Range("E$r").Value = $rowData.Hours
    $ws.Range("F$r").Value = $rowData.Type
    $ws.Range("G$r").Value = $rowData.Lecturer
    $ws.Range("H$r").Value = $rowData.Room
    $ws.Range("I$r").Value = $rowData.Block
    $ws.Range("J$r").Value = $rowData.Group
    $ws.Range("K$r").Value = $rowData.Level
    $ws.Range("L$r").Value = $rowData.Course
}
